$d = $word.ActiveDocument

function Replace-ParagraphXml($paragraphIndex, $xmlFragment) {
    $para = $d.Paragraphs.Item($paragraphIndex)
    $range = $para.Range
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $xmlFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- "Patcher" heading: lastRenderedPageBreak moves onto this bold run ---
Replace-ParagraphXml 66 '<w:p w:rsidR="00312CFE" w:rsidRPr="008B06FC" w:rsidRDefault="00312CFE" w:rsidP="00312CFE"><w:r w:rsidRPr="008B06FC"><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Patcher</w:t></w:r></w:p>'

# --- "Class function hooking" paragraph: lastRenderedPageBreak removed ---
Replace-ParagraphXml 67 '<w:p w:rsidR="00887B66" w:rsidRDefault="00887B66" w:rsidP="00887B66"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Class function hooking (ecx preservation).</w:t></w:r><w:r w:rsidR="005576D7"><w:t xml:space="preserve"> (xchg ecx, [esp]; push ecx)</w:t></w:r></w:p>'

# --- "Delay import directory." paragraph: lastRenderedPageBreak added ---
Replace-ParagraphXml 99 '<w:p w:rsidR="009420FA" w:rsidRPr="008B06FC" w:rsidRDefault="009420FA" w:rsidP="009420FA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="008B06FC"><w:lastRenderedPageBreak/><w:t>Delay import directory.</w:t></w:r></w:p>'

# --- "Bound import directory." paragraph: lastRenderedPageBreak removed ---
Replace-ParagraphXml 100 '<w:p w:rsidR="009420FA" w:rsidRPr="008B06FC" w:rsidRDefault="009420FA" w:rsidP="009420FA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="008B06FC"><w:t>Bound import directory.</w:t></w:r></w:p>'

# --- "Raise Intel warning level." paragraph: drop the _GoBack bookmark (it moves to the new paragraph below) ---
Replace-ParagraphXml 41 '<w:p w:rsidR="00686210" w:rsidRDefault="00686210" w:rsidP="00144F10"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Raise Intel warning level.</w:t></w:r></w:p>'

# --- Insert the two new list items after "Raise Intel warning level." ---
$p41 = $d.Paragraphs.Item(41)
$p41.Range.InsertParagraphAfter()

$p42 = $d.Paragraphs.Item(42)
$p42.Range.Text = "Stack trace on error."

$p42.Range.InsertParagraphAfter()

Replace-ParagraphXml 43 '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Improve process enumeration (improve ProcessEntry type</w:t></w:r><w:r><w:t xml:space="preserve"> and merge with ProcessList header</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>).</w:t></w:r></w:p>'

Write-Output "done"
